$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 288
$ws.Range("F3").Value = 256
$ws.Range("F4").Value = 293
$ws.Range("F5").Value = 2930
$ws.Range("F6").Value = 73
$ws.Range("F8").Value = 2287
$ws.Range("F9").Value = 1572
$ws.Range("F11").Value = 836
$ws.Range("F13").Value = 2620
$ws.Range("F15").Value = 1469
$ws.Range("F16").Value = 6861
$ws.Range("F18").Value = 7012
$ws.Range("F19").Value = 7012
$ws.Range("F21").Value = 2295
$ws.Range("F22").Value = 3053
$ws.Range("F23").Value = 3430
$ws.Range("F24").Value = 212
$ws.Range("F25").Value = 123
$ws.Range("F26").Value = 1770
$ws.Range("F28").Value = 288
$ws.Range("F29").Value = 863
$ws.Range("F30").Value = 6
$ws.Range("F32").Value = 30
$ws.Range("F33").Value = 373
$ws.Range("F34").Value = 1079
$ws.Range("F35").Value = 2436
$ws.Range("F36").Value = 13
$ws.Range("F37").Value = 154
$ws.Range("F38").Value = 349
$ws.Range("F39").Value = 962
$ws.Range("F40").Value = 197
$ws.Range("F41").Value = 445
$ws.Range("F42").Value = 494

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 67
$ws.Range("F21").Value = 59

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 67
$ws.Range("F4").Value = 288
$ws.Range("F7").Value = 293
$ws.Range("F8").Value = 2930
$ws.Range("F9").Value = 73
$ws.Range("F10").Value = 2287
$ws.Range("F11").Value = 1572
$ws.Range("F13").Value = 836
$ws.Range("F16").Value = 2620
$ws.Range("F17").Value = 1469
$ws.Range("F21").Value = 6861
$ws.Range("F23").Value = 7012
$ws.Range("F24").Value = 7012
$ws.Range("F26").Value = 2295
$ws.Range("F27").Value = 3053
$ws.Range("F28").Value = 3430
$ws.Range("F30").Value = 212
$ws.Range("F33").Value = 1770
$ws.Range("F36").Value = 288
$ws.Range("F37").Value = 863
$ws.Range("F39").Value = 30
$ws.Range("F40").Value = 373
$ws.Range("F41").Value = 59
$ws.Range("F42").Value = 2436
$ws.Range("F43").Value = 13
$ws.Range("F44").Value = 154
$ws.Range("F46").Value = 349
$ws.Range("F47").Value = 962
$ws.Range("F48").Value = 197
$ws.Range("F49").Value = 445
$ws.Range("F50").Value = 494
